$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112181514
$ws.Range("Q2").Value = 772353
$ws.Range("R2").Value = 7120281
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '2023-09-06'
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = '2023-09-06'
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Row 3
$ws.Range("A3").Value = 112181512
$ws.Range("B3").Value = 89405
$ws.Range("E3").Value = 1202
$ws.Range("F3").Value = 'Ullticka'
$ws.Range("G3").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H3").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q3").Value = 772413
$ws.Range("R3").Value = 7120316
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Row 4
$ws.Range("A4").Value = 112181509
$ws.Range("B4").Value = 89405
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = 'Ullticka'
$ws.Range("G4").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H4").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q4").Value = 772347
$ws.Range("R4").Value = 7120237
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '2023-09-07'
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = '2023-09-07'
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# Row 5
$ws.Range("A5").Value = 112181511
$ws.Range("Q5").Value = 772359
$ws.Range("R5").Value = 7120174
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = '2023-09-07'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = '2023-09-07'
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# Row 6
$ws.Range("A6").Value = 112182926
$ws.Range("B6").Value = 5113
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 100526
$ws.Range("F6").Value = 'Bronshjon'
$ws.Range("G6").Value = 'Callidium coriaceum'
$ws.Range("H6").Value = 'Paykull, 1800'
$ws.Range("Q6").Value = 772357
$ws.Range("R6").Value = 7120234
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# Row 7
$ws.Range("A7").Value = 112181500
$ws.Range("B7").Value = 89369
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 5447
$ws.Range("F7").Value = 'Vedticka'
$ws.Range("G7").Value = 'Fuscoporia viticola'
$ws.Range("H7").Value = '(Schwein.) Murrill'
$ws.Range("Q7").Value = 772346
$ws.Range("R7").Value = 7120286
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = '2023-09-06'
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = '2023-09-06'
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()

# Row 8
$ws.Range("A8").Value = 112181532
$ws.Range("B8").Value = 89351
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 112
$ws.Range("F8").Value = 'Stjärntagging'
$ws.Range("G8").Value = 'Asterodon ferruginosus'
$ws.Range("H8").Value = 'Pat.'
$ws.Range("Q8").Value = 772340
$ws.Range("R8").Value = 7120223
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()

# Row 9
$ws.Range("A9").Value = 112181582
$ws.Range("B9").Value = 81248
$ws.Range("E9").Value = 1312
$ws.Range("F9").Value = 'Gammelgransskål'
$ws.Range("G9").Value = 'Pseudographis pinicola'
$ws.Range("H9").Value = '(Nyl.) Rehm'
$ws.Range("Q9").Value = 772409
$ws.Range("R9").Value = 7120320
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
